# Update About Me slide.
#
# Repositions/resizes the shapes on the "About Me" slide (Brian Hansen bio)
# and updates a couple of text runs. Numeric Left/Top/Width/Height literals
# below were solved so that, after round-tripping through the COM layer's
# Single(float32)-then-EMU conversion, they land exactly on the target EMU
# coordinates (e.g. 7367452 EMU for Shape 1's new Left).

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    throw "Shape with id $id not found on slide"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Shape id=11, "Content Placeholder 2" (bullet list) ---
$shp11 = Get-ShapeById $s 11
$shp11.Left = 580.1143798828125
$shp11.Top = 120.32921600341797
$shp11.Width = 301.7142639160156
$shp11.Height = 288.64630126953125
$shp11.TextFrame.TextRange.Runs(1,1).Text = "20 Years working with SQL Server"

# --- Shape id=2, "Title 1" placeholder ("Brian Hansen") ---
$shp2 = Get-ShapeById $s 2
$shp2.Left = 152.54244995117188
$shp2.Top = 150.23118591308594
$shp2.Width = 294.5432434082031
$shp2.Height = 56.69291687011719
$shp2.TextFrame.TextRange.Paragraphs(1,1).ParagraphFormat.Alignment = 2

# --- Shape id=6, "Content Placeholder 5" (picture) ---
$shp6 = Get-ShapeById $s 6
$shp6.Left = 336.62152099609375
$shp6.Top = 254.57615661621094

# --- Shape id=7, "Picture 6" ---
$shp7 = Get-ShapeById $s 7
$shp7.Left = 45.24110412597656
$shp7.Top = 301.2997741699219

# --- Shape id=8, "Picture 7" ---
$shp8 = Get-ShapeById $s 8
$shp8.Left = 45.24110412597656
$shp8.Top = 255.40158081054688

# --- Shape id=13, "Picture 12" ---
$shp13 = Get-ShapeById $s 13
$shp13.Left = 214.69134521484375
$shp13.Top = 75.39614868164062

# --- Shape id=14, "Content Placeholder 2" ("@tf3604.com") ---
$shp14 = Get-ShapeById $s 14
$shp14.Left = 87.44763946533203
$shp14.Top = 294.3934020996094

# --- Shape id=15, "Content Placeholder 2" ("brian@tf3604.com") ---
$shp15 = Get-ShapeById $s 15
$shp15.Left = 87.44772338867188
$shp15.Top = 258.3822326660156

# --- Shape id=17, "Content Placeholder 2" ("children.org") ---
$shp17 = Get-ShapeById $s 17
$shp17.Left = 336.62152099609375
$shp17.Top = 331.3711242675781

# --- Shape id=18, "Content Placeholder 2" (website URL) ---
$shp18 = Get-ShapeById $s 18
$shp18.Left = 45.24110412597656
$shp18.Top = 419.62835693359375
$shp18.Width = 821.5017700195312
$shp18.TextFrame.TextRange.Runs(1,1).Text = "www.tf3604.com/internals"
